$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Header row additions / changes
$ws.Range("D1").Value = "COMMENT"

# Fill COMMENT column for all data rows (existing rows 2-8 + new rows 9-15)
$ws.Range("D2").Value = "NO COMMENT"
$ws.Range("D3").Value = "NO COMMENT"
$ws.Range("D4").Value = "NO COMMENT"
$ws.Range("D5").Value = "NO COMMENT"
$ws.Range("D6").Value = "NO COMMENT"
$ws.Range("D7").Value = "NO COMMENT"
$ws.Range("D8").Value = "NO COMMENT"
$ws.Range("D9").Value = "NO COMMENT"
$ws.Range("D10").Value = "NO COMMENT"
$ws.Range("D11").Value = "NO COMMENT"
$ws.Range("D12").Value = "NO COMMENT"
$ws.Range("D13").Value = "NO COMMENT"
$ws.Range("D14").Value = "NO COMMENT"
$ws.Range("D15").Value = "NO COMMENT"

# New AF102 device rows (9-15): DEVICE column
$ws.Range("A9").Value = "AF102"
$ws.Range("A10").Value = "AF102"
$ws.Range("A11").Value = "AF102"
$ws.Range("A12").Value = "AF102"
$ws.Range("A13").Value = "AF102"
$ws.Range("A14").Value = "AF102"
$ws.Range("A15").Value = "AF102"

# New AF102 device rows (9-15): NAME column
$ws.Range("C9").Value = "TATU"
$ws.Range("C10").Value = "SAMI"
$ws.Range("C11").Value = "TEEMU"
$ws.Range("C12").Value = "VICE"
$ws.Range("C13").Value = "SMILEY"
$ws.Range("C14").Value = "TIPZU"
$ws.Range("C15").Value = "FINU"

# IOCARD column: DI -> BOOL for all rows 2-15
$ws.Range("B2").Value = "BOOL"
$ws.Range("B3").Value = "BOOL"
$ws.Range("B4").Value = "BOOL"
$ws.Range("B5").Value = "BOOL"
$ws.Range("B6").Value = "BOOL"
$ws.Range("B7").Value = "BOOL"
$ws.Range("B8").Value = "BOOL"
$ws.Range("B9").Value = "BOOL"
$ws.Range("B10").Value = "BOOL"
$ws.Range("B11").Value = "BOOL"
$ws.Range("B12").Value = "BOOL"
$ws.Range("B13").Value = "BOOL"
$ws.Range("B14").Value = "BOOL"
$ws.Range("B15").Value = "BOOL"

# Column C width
$ws.Columns.Item(3).ColumnWidth = 18.5703125

# Activate Sheet1 and set selection
$ws.Activate() | Out-Null
$ws.Range("F14:G14").Select() | Out-Null
